# "Add mobile handson 2nd"
# Slide 3 ("TIS株式会社のインターンシップに参加しませんか？") gains an extra
# line explaining the event is online, and the text box grows to fit it.

$p = $ppt.ActivePresentation
$s = $p.Slides.Item(3)
$sh = $s.Shapes.Item(3)

$tr = $sh.TextFrame.TextRange
$firstPara = $tr.Paragraphs(1, 1)

# Split right after the first paragraph ("TIS株式会社のインターンシップに
# 参加しませんか？") so a brand-new paragraph is created after it.
$firstPara.InsertAfter("`r") | Out-Null

# Fill the freshly created (now second) paragraph with the new sentence and
# match the bold, 16pt, Meiryo UI formatting used by the first paragraph.
$secondPara = $tr.Paragraphs(2, 1)
$newRun = $secondPara.InsertBefore("オンライン開催なので会津から参加できます。")
$newRun.Font.Bold = $true

# Grow the text box to fit the extra line (EMU -> points: /12700). Set this
# after the text edit since the autosize box recomputes height on change.
$sh.Height = 6185989 / 12700
